# Apply the "Add files via upload / finalssasada" edit:
#  - Update longitude (B) and latitude (C) values for rows 2-11
#  - Re-type the county names in column G using Title Case (was lower-case)
#  - Apply a custom 4-decimal number format to B11
#  - Resize column B and change the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated longitude / latitude values (column B / C), rows 2-11 ---
$ws.Range("B2").Value = -119.456
$ws.Range("C2").Value = 34.429040000000001

$ws.Range("B3").Value = -121.9426
$ws.Range("C3").Value = 39.17313

$ws.Range("B4").Value = -121.4842
$ws.Range("C4").Value = 37.448520000000002

$ws.Range("B5").Value = -122.3036
$ws.Range("C5").Value = 37.144849999999998

$ws.Range("B6").Value = -121.5874
$ws.Range("C6").Value = 38.004289999999997

$ws.Range("B7").Value = -121.83369999999999
$ws.Range("C7").Value = 38.061590000000002

$ws.Range("B8").Value = -121.2551
$ws.Range("C8").Value = 37.024529999999999

$ws.Range("B9").Value = -120.5331
$ws.Range("C9").Value = 37.035989999999998

$ws.Range("B10").Value = -120.29819999999999
$ws.Range("C10").Value = 35.901539999999997

$ws.Range("B11").Value = -120.8998
$ws.Range("C11").Value = 36.738050000000001

# --- Re-typed (Title Case) county / subregion names, column G ---
$ws.Range("G2").Value = "Ventura"
$ws.Range("G3").Value = "Sutter"
$ws.Range("G4").Value = "Stanislaus"
$ws.Range("G5").Value = "Santa Cruz"
$ws.Range("G6").Value = "San Joaquin"
$ws.Range("G7").Value = "Sacramento"
$ws.Range("G8").Value = "Merced"
$ws.Range("G9").Value = "Madera"
$ws.Range("G10").Value = "Kings"
$ws.Range("G11").Value = "Fresno"

# --- Custom number format (4 decimal places) applied to B11 ---
$ws.Range("B11").NumberFormat = "0.0000"

# --- Column B width ---
$ws.Columns.Item(2).ColumnWidth = 11.1667

# --- Update the selected cell ---
$ws.Range("C2").Select() | Out-Null
